$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.543.31'
$ws.Range('E2').Value = '  +3.17%  '
$ws.Range('D3').Value = '1.844.69'
$ws.Range('E3').Value = '  +2.47%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.80'
$ws.Range('E5').Value = '  +3.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.619'
$ws.Range('E6').Value = '  +2.97%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.62'
$ws.Range('E8').Value = '  +14.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.312'
$ws.Range('E9').Value = '  +8.38%  '
$ws.Range('E10').Value = '  +4.70%  '
$ws.Range('E11').Value = '  +2.71%  '
$ws.Range('D12').Value = '2.111.33'
$ws.Range('E12').Value = '  +2.48%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.849.38'
$ws.Range('E13').Value = '  +2.68%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.34'
$ws.Range('E14').Value = '  +4.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.676'
$ws.Range('E15').Value = '  +7.73%  '
$ws.Range('E16').Value = '  +8.61%  '
$ws.Range('D17').Value = '35.529.75'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.60'
$ws.Range('E18').Value = '  +3.93%  '
$ws.Range('E19').Value = '  +5.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '244.58'
$ws.Range('E20').Value = '  +2.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.11'
$ws.Range('E21').Value = '  +9.36%  '
$ws.Range('E22').Value = '  +13.77%  '
$ws.Range('E24').Value = '  +5.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.29'
$ws.Range('E25').Value = '  +0.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.03'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.86'
$ws.Range('E27').Value = '  +1.84%  '
$ws.Range('E28').Value = '  +1.28%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.57'
$ws.Range('E29').Value = '  +28.38%  '
$ws.Range('E30').Value = '  +0.12%  '
$ws.Range('D31').Value = '3.319.58'
$ws.Range('E31').Value = '  +36.63%  '
$ws.Range('E32').Value = '  +8.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.14'
$ws.Range('E33').Value = '  +8.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.95'
$ws.Range('E34').Value = '  +5.74%  '
$ws.Range('E35').Value = '  +2.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '95.87'
$ws.Range('E36').Value = '  +17.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.696'
$ws.Range('E37').Value = '  +9.20%  '
$ws.Range('E38').Value = '  +8.74%  '
$ws.Range('D39').Value = '1.350.72'
$ws.Range('E39').Value = '  +3.41%  '
$ws.Range('E40').Value = '  +5.75%  '
$ws.Range('E41').Value = '  +6.35%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.02'
$ws.Range('E42').Value = '  +8.05%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '15.32'
$ws.Range('E43').Value = '  +9.22%  '
$ws.Range('E44').Value = '  +3.27%  '
$ws.Range('E45').Value = '  +0.78%  '
$ws.Range('E46').Value = '  +0.23%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0519'
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('D49').Value = '2.016.00'
$ws.Range('E49').Value = '  +2.86%  '
$ws.Range('E50').Value = '  +0.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '102.63'
$ws.Range('E51').Value = '  +0.88%  '
